# MVP for adding more pounds to specialty nforx
# Grapes (row 24) gets the "nforx" (n-for-x) specialty deal turned on,
# matching the pattern already used by Bacon (row 11): Limit=6,
# Specialty Variable 1=3, Specialty Variable 2=8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("F24").Value = $true      # Has Specialty -> TRUE
$ws.Range("G24").Value = "nforx"    # Type of Specialty
$ws.Range("H24").Value = 6          # Limit
$ws.Range("I24").Value = 3          # Specialty Variable 1
$ws.Range("J24").Value = 8          # Specialty Variable 2

# Update the view state: scroll so C7 is the top-left visible cell, then
# land the active selection on K24 (matches the author's recorded cursor
# position when they made this edit).
$ws.Range("C7").Select() | Out-Null
$ws.Range("K24").Select() | Out-Null
